# Updated cryptos list on Wed Jul 31 06:18:28 UTC 2024 with GitHub Actions
# Refreshes price/volume columns (D/E) for each coin row and reorders the
# Fetch.AI / Aptos rows (35/36) to match the latest ranking snapshot.
#
# Note: several "Price" values are plain numerals (e.g. "586.67", "0.999")
# that Excel would otherwise auto-convert to a Number on assignment. Those
# are written with a leading apostrophe (forces text entry, like typing it
# in the UI) and then the cell style is reset to "Normal" so the sheet ends
# up with the same formatting as before - only the literal text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.392.25"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "3.313.18"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'586.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "

$ws.Range("D6").Value = "'182.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.641"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.12%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.126"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.18%  "

$ws.Range("D10").Value = "'6.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("D11").Value = "'0.404"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").Value = "3.888.24"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("D14").Value = "66.443.71"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").Value = "'26.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.60%  "

$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "3.261.12"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").Value = "'431.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "'13.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("D20").Value = "'5.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.57%  "

$ws.Range("D21").Value = "'7.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.67%  "

$ws.Range("D22").Value = "'72.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'5.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("D25").Value = "3.432.01"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").Value = "'0.516"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("E27").Value = "  +2.34%  "

$ws.Range("D28").Value = "'0.0000114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.58%  "

$ws.Range("D29").Value = "'9.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").Value = "'22.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'5.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'6.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").Value = "'159.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "

$ws.Range("D39").Value = "'1.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").Value = "'26.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("D41").Value = "2.868.57"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").Value = "'0.772"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("D43").Value = "'4.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("D44").Value = "'40.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").Value = "'0.0667"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("D46").Value = "'6.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "

$ws.Range("E47").Value = "  -1.92%  "

$ws.Range("D48").Value = "'23.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("D49").Value = "'318.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").Value = "'0.0271"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("E51").Value = "  +3.67%  "
